$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hours for Matthew Handley (row 4) and Kathryn Swineford (row 6)
$ws.Range("C4").Value = 3.25
$ws.Range("C6").Value = 3.25

# Update the selected cell to match the saved view state
$ws.Range("F20").Select()
